# Manual continued. Chapter API completed as draft
# Append two new effort-log entries (rows 28 and 29) to the "effort" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: 2012-10-18, 1.75h, "Manual continued"
$ws.Cells.Item(28, 1).Value = 41200
$ws.Cells.Item(28, 2).Value = 1.75
$ws.Cells.Item(28, 4).Value = "Manual continued"

# Row 29: 2012-10-23, 2.5h, "Manual continued"
$ws.Cells.Item(29, 1).Value = 41205
$ws.Cells.Item(29, 2).Value = 2.5
$ws.Cells.Item(29, 4).Value = "Manual continued"

# Reflect the cursor ending up on the last edited cell, as in the authored file.
$ws.Range("D29").Select() | Out-Null
